$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 139.15384
$ws.Range("I9").Value = 137.09091
$ws.Range("K9").Value = 137.09091
$ws.Range("M9").Value = 31.90908999999999
$ws.Range("H17").Value = 34802
$ws.Range("I17").Value = 1200
$ws.Range("J17").Value = 36094.383
$ws.Range("K17").Value = 3600
$ws.Range("L17").Value = 108283.149
$ws.Range("M17").Value = -3432
$ws.Range("N17").Value = -108619.149
$ws.Range("H70").Value = 2466
$ws.Range("J70").Value = 2466
$ws.Range("L70").Value = 7398
$ws.Range("N70").Value = -7938
$ws.Range("H73").Value = 2466
$ws.Range("J73").Value = 2466
$ws.Range("L73").Value = 7398
$ws.Range("N73").Value = -9270
$ws.Range("H103").Value = 33333826
$ws.Range("J103").Value = 35714780
$ws.Range("L103").Value = 107144340
$ws.Range("N103").Value = -107145512
$ws.Range("H132").Value = 1936.4286
$ws.Range("I132").Value = 1253.5217
$ws.Range("J132").Value = 5077.8
$ws.Range("K132").Value = 3760.5651
$ws.Range("L132").Value = 15233.4
$ws.Range("M132").Value = -1230.5651
$ws.Range("N132").Value = -20293.4
$ws.Range("H135").Value = 1158.7307
$ws.Range("I135").Value = 954.8946999999999
$ws.Range("J135").Value = 1712
$ws.Range("K135").Value = 8594.052299999999
$ws.Range("L135").Value = 15408
$ws.Range("M135").Value = -6059.052299999999
$ws.Range("N135").Value = -20478
$ws.Range("H138").Value = 2391.4
$ws.Range("I138").Value = 1283.9474
$ws.Range("J138").Value = 3200.6924
$ws.Range("K138").Value = 3851.8422
$ws.Range("L138").Value = 9602.0772
$ws.Range("M138").Value = 1288.1578
$ws.Range("N138").Value = -19882.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 1499
$ws.Range("I23").Value = 1499
$ws.Range("K23").Value = 1499
$ws.Range("M23").Value = -1240
$ws.Range("H32").Value = 63218.43
$ws.Range("I32").Value = 44037.957
$ws.Range("J32").Value = 105066.73
$ws.Range("K32").Value = 44037.957
$ws.Range("L32").Value = 105066.73
$ws.Range("M32").Value = -43750.957
$ws.Range("N32").Value = -105640.73
$ws.Range("H61").Value = 1775.8928
$ws.Range("I61").Value = 1666.3462
$ws.Range("J61").Value = 3200
$ws.Range("K61").Value = 1666.3462
$ws.Range("L61").Value = 3200
$ws.Range("M61").Value = -1454.3462
$ws.Range("N61").Value = -3624
$ws.Range("H74").Value = 2253.8635
$ws.Range("I74").Value = 2123.1428
$ws.Range("J74").Value = 4999
$ws.Range("K74").Value = 2123.1428
$ws.Range("L74").Value = 4999
$ws.Range("M74").Value = -1249.1428
$ws.Range("N74").Value = -6747
$ws.Range("H77").Value = 2253.8635
$ws.Range("I77").Value = 2123.1428
$ws.Range("J77").Value = 4999
$ws.Range("K77").Value = 10615.714
$ws.Range("L77").Value = 24995
$ws.Range("M77").Value = -6247.714
$ws.Range("N77").Value = -33731
$ws.Range("H92").Value = 64275
$ws.Range("J92").Value = 64275
$ws.Range("L92").Value = 64275
$ws.Range("N92").Value = -69267
$ws.Range("H96").Value = 25344
$ws.Range("J96").Value = 25344
$ws.Range("L96").Value = 25344
$ws.Range("N96").Value = -30836
$ws.Range("H126").Value = 5016.25
$ws.Range("I126").Value = 5016.25
$ws.Range("K126").Value = 15048.75
$ws.Range("M126").Value = -12578.75
$ws.Range("H136").Value = 1775.8928
$ws.Range("I136").Value = 1666.3462
$ws.Range("J136").Value = 3200
$ws.Range("K136").Value = 4999.0386
$ws.Range("L136").Value = 9600
$ws.Range("M136").Value = -2449.0386
$ws.Range("N136").Value = -14700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 900000000
$ws.Range("J19").Value = 900000000
$ws.Range("L19").Value = 900000000
$ws.Range("N19").Value = -900000346
$ws.Range("H94").Value = 562.0323
$ws.Range("I94").Value = 600.10345
$ws.Range("J94").Value = 10
$ws.Range("K94").Value = 600.10345
$ws.Range("L94").Value = 10
$ws.Range("M94").Value = -149.10345
$ws.Range("N94").Value = -912
$ws.Range("H95").Value = 31978
$ws.Range("J95").Value = 31978
$ws.Range("L95").Value = 31978
$ws.Range("N95").Value = -37470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 35000000
$ws.Range("J4").Value = 35000000
$ws.Range("L4").Value = 35000000
$ws.Range("N4").Value = -35000224
$ws.Range("H15").Value = 550
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 550
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = $null
$ws.Range("M15").Value = 550
$ws.Range("N15").Value = -890
$ws.Range("H31").Value = 2682.8333
$ws.Range("I31").Value = 2016.3636
$ws.Range("K31").Value = 2016.3636
$ws.Range("M31").Value = -1721.3636
$ws.Range("H34").Value = 2682.8333
$ws.Range("I34").Value = 2016.3636
$ws.Range("K34").Value = 2016.3636
$ws.Range("M34").Value = -1814.3636
$ws.Range("H62").Value = 4013.7273
$ws.Range("I62").Value = 4343.5713
$ws.Range("J62").Value = 3436.5
$ws.Range("K62").Value = 4343.5713
$ws.Range("L62").Value = 3436.5
$ws.Range("M62").Value = -3719.5713
$ws.Range("N62").Value = -4684.5
$ws.Range("H65").Value = 4013.7273
$ws.Range("I65").Value = 4343.5713
$ws.Range("J65").Value = 3436.5
$ws.Range("K65").Value = 21717.8565
$ws.Range("L65").Value = 17182.5
$ws.Range("M65").Value = -18597.8565
$ws.Range("N65").Value = -23422.5
$ws.Range("H132").Value = 3284
$ws.Range("I132").Value = 3284
$ws.Range("K132").Value = 9852
$ws.Range("M132").Value = -7322

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4936715.5
$ws.Range("I4").Value = 1259555.5
$ws.Range("K4").Value = 3778666.5
$ws.Range("M4").Value = -3778554.5
$ws.Range("H38").Value = 75.666664
$ws.Range("I38").Value = 75.333336
$ws.Range("J38").Value = 76.333336
$ws.Range("K38").Value = 226.000008
$ws.Range("L38").Value = 229.000008
$ws.Range("M38").Value = 120.999992
$ws.Range("N38").Value = -923.000008
$ws.Range("H68").Value = 1456.8572
$ws.Range("I68").Value = 1574.5
$ws.Range("K68").Value = 4723.5
$ws.Range("M68").Value = -3912.5
$ws.Range("H69").Value = 5770.3335
$ws.Range("J69").Value = 6500.5
$ws.Range("L69").Value = 19501.5
$ws.Range("N69").Value = -21123.5
$ws.Range("H71").Value = 1456.8572
$ws.Range("I71").Value = 1574.5
$ws.Range("K71").Value = 14170.5
$ws.Range("M71").Value = -10114.5
$ws.Range("H72").Value = 5770.3335
$ws.Range("J72").Value = 6500.5
$ws.Range("L72").Value = 58504.5
$ws.Range("N72").Value = -66616.5
$ws.Range("H132").Value = 2597.457
$ws.Range("J132").Value = 2873.12
$ws.Range("L132").Value = 25858.08
$ws.Range("N132").Value = -30918.08

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = $null
$ws.Range("H80").Value = 4774.5557
$ws.Range("I80").Value = 3996.8333
$ws.Range("K80").Value = 3996.8333
$ws.Range("M80").Value = -2998.8333
$ws.Range("H83").Value = 4774.5557
$ws.Range("I83").Value = 3996.8333
$ws.Range("K83").Value = 19984.1665
$ws.Range("M83").Value = -14992.1665
$ws.Range("H102").Value = 3580.2727
$ws.Range("I102").Value = 3580.2727
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3580.2727
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = $null
$ws.Range("N102").Value = -1958.2727
$ws.Range("H113").Value = 3860.2
$ws.Range("J113").Value = 5995
$ws.Range("L113").Value = 5995
$ws.Range("N113").Value = -10335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1212.5
$ws.Range("H27").Value = 1212.5
$ws.Range("H55").Value = 998.2857
$ws.Range("J55").Value = 1247
$ws.Range("L55").Value = 1247
$ws.Range("N55").Value = -1593
$ws.Range("H122").Value = 12257.833
$ws.Range("I122").Value = 14953.889
$ws.Range("K122").Value = 44861.667
$ws.Range("M122").Value = -42411.667
$ws.Range("H136").Value = 5013.0713
$ws.Range("J136").Value = 5609.4443
$ws.Range("L136").Value = 16828.3329
$ws.Range("N136").Value = -21928.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 4222675.5
$ws.Range("J3").Value = 26000
$ws.Range("L3").Value = 26000
$ws.Range("N3").Value = -26228
$ws.Range("H14").Value = 1002999.2
$ws.Range("I14").Value = 6000000
$ws.Range("J14").Value = 3599
$ws.Range("K14").Value = 6000000
$ws.Range("L14").Value = 3599
$ws.Range("M14").Value = -5999832
$ws.Range("N14").Value = -3935
$ws.Range("H23").Value = 2950
$ws.Range("I23").Value = 2950
$ws.Range("K23").Value = 2950
$ws.Range("M23").Value = -2721
$ws.Range("H41").Value = 19690.3
$ws.Range("J41").Value = 20980.857
$ws.Range("L41").Value = 20980.857
$ws.Range("N41").Value = -21760.857
$ws.Range("H81").Value = 3942.4666
$ws.Range("I81").Value = 4325.2173
$ws.Range("J81").Value = 2684.8572
$ws.Range("K81").Value = 8650.434600000001
$ws.Range("L81").Value = 5369.7144
$ws.Range("M81").Value = -7589.434600000001
$ws.Range("N81").Value = -7491.7144
$ws.Range("H84").Value = 3942.4666
$ws.Range("I84").Value = 4325.2173
$ws.Range("J84").Value = 2684.8572
$ws.Range("K84").Value = 43252.173
$ws.Range("L84").Value = 26848.572
$ws.Range("M84").Value = -37948.173
$ws.Range("N84").Value = -37456.572
$ws.Range("H122").Value = 1687.8948
$ws.Range("I122").Value = 1726.7858
$ws.Range("J122").Value = 1579
$ws.Range("K122").Value = 5180.357400000001
$ws.Range("L122").Value = 4737
$ws.Range("M122").Value = -2730.357400000001
$ws.Range("N122").Value = -9637
